$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kandam2")
Write-Host $ws.Cells.Item(44, 3).Value
